$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseRun")

# Pre-format the new data range as Text so numeric-looking IDs
# (e.g. "11539914") are stored as strings, not numbers.
$ws.Range("A6:E17").NumberFormat = "@"

$ws.Range("A6").Value = "TN2485352"
$ws.Range("B6").Value = "11539914"
$ws.Range("C6").Value = "Personal Auto - Credit"
$ws.Range("D6").Value = "Base"
$ws.Range("E6").Value = "TC002"

$ws.Range("A7").Value = "TN2485356"
$ws.Range("B7").Value = "11539933"
$ws.Range("C7").Value = "Personal Auto - Credit"
$ws.Range("D7").Value = "Restricted"
$ws.Range("E7").Value = "TC003"

$ws.Range("A8").Value = "TB2485358"
$ws.Range("B8").Value = "11539965"
$ws.Range("C8").Value = "Bond - No Credit"
$ws.Range("D8").Value = "Nil"
$ws.Range("E8").Value = "TC001"

$ws.Range("A9").Value = "TN2485359"
$ws.Range("B9").Value = "11539970"
$ws.Range("C9").Value = "Personal Auto - Credit"
$ws.Range("D9").Value = "Base"
$ws.Range("E9").Value = "TC002"

$ws.Range("A10").Value = "TN2485361"
$ws.Range("B10").Value = "11539983"
$ws.Range("C10").Value = "Personal Auto - Credit"
$ws.Range("D10").Value = "Restricted"
$ws.Range("E10").Value = "TC003"

$ws.Range("A11").Value = "TN2485363"
$ws.Range("B11").Value = "11540001"
$ws.Range("C11").Value = "Personal Auto - Credit"
$ws.Range("D11").Value = "Restricted"
$ws.Range("E11").Value = "TC003"

$ws.Range("A12").Value = "TN2485392"
$ws.Range("B12").Value = "11540661"
$ws.Range("C12").Value = "Personal Auto - Credit"
$ws.Range("D12").Value = "Restricted"
$ws.Range("E12").Value = "TC003"

$ws.Range("A13").Value = "TN2485397"
$ws.Range("B13").Value = "11540740"
$ws.Range("C13").Value = "Personal Auto - Credit"
$ws.Range("D13").Value = "Base"
$ws.Range("E13").Value = "TC002"

$ws.Range("A14").Value = "TN2485399"
$ws.Range("B14").Value = "11540748"
$ws.Range("C14").Value = "Personal Auto - Credit"
$ws.Range("D14").Value = "Base"
$ws.Range("E14").Value = "TC002"

$ws.Range("A15").Value = "TN2485403"
$ws.Range("B15").Value = "11540787"
$ws.Range("C15").Value = "Personal Auto - Credit"
$ws.Range("D15").Value = "Base"
$ws.Range("E15").Value = "TC002"

$ws.Range("A16").Value = "TN2485405"
$ws.Range("B16").Value = "11540797"
$ws.Range("C16").Value = "Personal Auto - Credit"
$ws.Range("D16").Value = "Base"
$ws.Range("E16").Value = "TC002"

$ws.Range("A17").Value = "TN2485407"
$ws.Range("B17").Value = "11540814"
$ws.Range("C17").Value = "Personal Auto - Credit"
$ws.Range("D17").Value = "Base"
$ws.Range("E17").Value = "TC002"

# Reset the style back to Normal/General so no extra cell styling
# is left behind (only the underlying value stays textual).
$ws.Range("A6:E17").Style = "Normal"
